# Words1.xlsx - fix spelling mistakes in the English column of the word list,
# then leave the selection on the cell that was last corrected (C14),
# matching the author's editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "centers" -> "centres" (row 14)
$ws.Range("C14").Value = "centres"

# "olil" -> "oil" (row 22)
$ws.Range("C22").Value = "oil"

# "to notifiy" -> "to notify" (row 26)
$ws.Range("C26").Value = "to notify"

# Leave the active selection where the edits were last made.
$ws.Range("C14").Select()
